$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5-27 down to 6-28
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new record
$ws.Cells.Item(5, 1).Value = 5
$ws.Cells.Item(5, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(5, 3).Value = "Maule"
$ws.Cells.Item(5, 4).Value = 44635
$ws.Cells.Item(5, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(5, 5).Value = 7
$ws.Cells.Item(5, 6).Value = 100112043
$ws.Cells.Item(5, 7).Value = "Pepino dulce"
$ws.Cells.Item(5, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 300
$ws.Cells.Item(5, 11).Value = 15000
$ws.Cells.Item(5, 12).Value = 15000
$ws.Cells.Item(5, 13).Value = 15000
$ws.Cells.Item(5, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(5, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(5, 16).Value = 833
$ws.Cells.Item(5, 17).Value = 18
$ws.Cells.Item(5, 18).Value = "Hortaliza"
